$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new "0.1.7" version row (row 10) plus a fresh blank row (row 11) ---
# Copy formatting (fill/border/alignment/number-format) from the existing
# even-row template (row 8 -> A-D style 4, E-G style 8) down into row 10,
# and from the odd-row template (row 9 -> A-D style 2, E-G style 3) into the
# new trailing blank row 11 - this mirrors how every prior version entry in
# this sheet was appended (alternating banding, same column styles).
$ws.Range("A8:G8").Copy()
$ws.Range("A10:G10").PasteSpecial(-4122)
$ws.Range("A9:G9").Copy()
$ws.Range("A11:G11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 10 values ---
# (The shared-string table is append-only in authoring order, so D10 is
# written before C10 here to land the two new paragraph strings at indices
# 43/44 in the same order the source workbook uses.)
$ws.Range("A10").Value = "0.1.7"
$ws.Range("B10").Value = "AUTOMATA CELULAR - copia (14)"
$ws.Range("D10").Value = "-Association to itself corrected.`n-Aggrupation to itself corrected.`n-SG corrected to random order when they have the same value.`n-Added Agrupation and Reciprocal to SG.`n-Reciprocal working.`n-UI: automatically fill aggrupation data."
$ws.Range("C10").Value = "-Agrupation and disaggregate to be done in functions.`n-Change reproduction and distribution to two parts.`n-UI: Delete rows according to working functionality.`n*Implement mutations.`n-With 4 or less niches the distribution is not equaly done.`n-Document every function.`n-Disaggregate to be done in between reproduction and distribution."
$ws.Range("E10").Value = $ws.Range("E9").Value2
$ws.Range("F10").Value = $ws.Range("F9").Value2
$ws.Range("G10").Value = $ws.Range("G9").Value2

# Row height for the new data row mirrors the auto-fit height Excel applies
# to wrapped, multi-line cells (matches the 100.8pt the real sheet reaches).
$ws.Range("A10:G10").RowHeight = 100.8

# Row 11 stays an empty styled row (same pattern as the sheet's very first
# templated row 2), ready for the next version entry.

# --- Update the active selection to follow the newly appended rows ---
$ws.Range("C14").Select()
